{"js": "// Redefine the Synopsis front matter:\n//  1. Split the \"ACKNOWLEDGEMENT\" heading paragraph into a leading empty\n//     paragraph (keeping the paragraph-mark formatting) plus the heading\n//     paragraph itself, which now lands at the top of a rendered page and\n//     therefore carries a <w:lastRenderedPageBreak/> marker.\n//  2. The \"Title: Sign Language Translator using AI \" paragraph no longer\n//     starts a rendered page, so its <w:lastRenderedPageBreak/> marker is\n//     removed.\n//\n// Each step re-reads the paragraph collection right before acting on it so\n// that the second edit never operates on a paragraph reference that was\n// invalidated by the first (structure-changing) edit.\n\nfunction findParagraphByExactText(items, text) {\n  for (const paragraph of items) {\n    if (paragraph.text === text) return paragraph;\n  }\n  return null;\n}\n\nfunction findParagraphByPrefix(items, prefix) {\n  for (const paragraph of items) {\n    if (paragraph.text.indexOf(prefix) === 0) return paragraph;\n  }\n  return null;\n}\n\nconst body = context.document.body;\n\n// ------------------------------------------------------------------\n// 1) Split the ACKNOWLEDGEMENT heading paragraph in two.\n// ------------------------------------------------------------------\nlet paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst ackParagraph = findParagraphByExactText(paragraphs.items, \"ACKNOWLEDGEMENT\");\n\nif (ackParagraph) {\n  const ackRange = ackParagraph.getRange(\"Whole\");\n\n  const ackOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Heading3\"/>\n              <w:spacing w:before=\"47\"/>\n              <w:rPr>\n                <w:color w:val=\"2E5395\"/>\n              </w:rPr>\n            </w:pPr>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Heading3\"/>\n              <w:spacing w:before=\"47\"/>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:color w:val=\"2E5395\"/>\n              </w:rPr>\n              <w:lastRenderedPageBreak/>\n              <w:t>ACKNOWLEDGEMENT</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\n  ackRange.insertOoxml(ackOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ------------------------------------------------------------------\n// 2) Drop the stale <w:lastRenderedPageBreak/> on the synopsis title.\n//    Re-fetch the paragraph collection since edit (1) may have shifted it.\n// ------------------------------------------------------------------\nparagraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst titleParagraph = findParagraphByPrefix(\n  paragraphs.items,\n  \"Title: Sign Language Translator using AI\"\n);\n\nif (titleParagraph) {\n  const titleRange = titleParagraph.getRange(\"Whole\");\n\n  const titleOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:jc w:val=\"center\"/>\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n                <w:sz w:val=\"32\"/>\n                <w:szCs w:val=\"32\"/>\n                <w:u w:val=\"single\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n                <w:sz w:val=\"36\"/>\n                <w:szCs w:val=\"36\"/>\n                <w:u w:val=\"single\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">Title: Sign Language Translator using AI </w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\n  titleRange.insertOoxml(titleOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# ------------------------------------------------------------------\n# 1) Split the \"ACKNOWLEDGEMENT\" heading paragraph into two paragraphs:\n#    an empty paragraph carrying the paragraph-mark run formatting,\n#    followed by the heading paragraph itself (now starting a new\n#    rendered page, hence the <w:lastRenderedPageBreak/>).\n# ------------------------------------------------------------------\n$ackPara = $null\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text -eq \"ACKNOWLEDGEMENT`r\") {\n        $ackPara = $p\n        break\n    }\n}\n\nif ($ackPara -ne $null) {\n    $ackRange = $ackPara.Range\n\n    $ackXml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Heading3\"/>\n              <w:spacing w:before=\"47\"/>\n              <w:rPr>\n                <w:color w:val=\"2E5395\"/>\n              </w:rPr>\n            </w:pPr>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Heading3\"/>\n              <w:spacing w:before=\"47\"/>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:color w:val=\"2E5395\"/>\n              </w:rPr>\n              <w:lastRenderedPageBreak/>\n              <w:t>ACKNOWLEDGEMENT</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n    $ackRange.InsertXML($ackXml)\n}\n\n# ------------------------------------------------------------------\n# 2) The synopsis title paragraph no longer starts a rendered page,\n#    so drop its <w:lastRenderedPageBreak/> marker.\n# ------------------------------------------------------------------\n$titlePara = $null\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text -eq \"Title: Sign Language Translator using AI `r\") {\n        $titlePara = $p\n        break\n    }\n}\n\nif ($titlePara -ne $null) {\n    $titleRange = $titlePara.Range\n\n    $titleXml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:jc w:val=\"center\"/>\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n                <w:sz w:val=\"32\"/>\n                <w:szCs w:val=\"32\"/>\n                <w:u w:val=\"single\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:b/>\n                <w:bCs/>\n                <w:sz w:val=\"36\"/>\n                <w:szCs w:val=\"36\"/>\n                <w:u w:val=\"single\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">Title: Sign Language Translator using AI </w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n    $titleRange.InsertXML($titleXml)\n}\n"}
